$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "historico_respostas" rows appended after the previous last row (203),
# i.e. sequential ids 217..236 with the matching date/questao/disciplina/
# assunto/resultado values from the commit.
$data = @(
    @(217, "10/01/2026 04:39", 871, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(218, "10/01/2026 04:41", 794, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(219, "10/01/2026 04:43", 504, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(220, "10/01/2026 04:45", 818, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(221, "10/01/2026 04:46", 499, "Conhecimentos Específicos", "Engenharia Organizacional", 0),
    @(222, "10/01/2026 04:49", 816, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(223, "10/01/2026 04:51", 500, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(224, "10/01/2026 04:52", 779, "Conhecimentos Específicos", "Engenharia Organizacional", 0),
    @(225, "10/01/2026 04:53", 870, "Conhecimentos Específicos", "Engenharia Organizacional", 0),
    @(226, "10/01/2026 04:57", 505, "Conhecimentos Específicos", "Engenharia Organizacional", 1),
    @(227, "10/01/2026 05:28", 856, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(228, "10/01/2026 05:37", 868, "Conhecimentos Específicos", "Gestão de Estoques", 1),
    @(229, "10/01/2026 05:51", 774, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(230, "10/01/2026 05:53", 447, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(231, "10/01/2026 05:59", 801, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(232, "10/01/2026 06:10", 449, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(233, "10/01/2026 06:18", 768, "Conhecimentos Específicos", "Gestão de Estoques", 1),
    @(234, "10/01/2026 06:19", 832, "Conhecimentos Específicos", "Gestão de Estoques", 0),
    @(235, "10/01/2026 06:20", 450, "Conhecimentos Específicos", "Gestão de Estoques", 1),
    @(236, "10/01/2026 06:20", 1242, "Conhecimentos Específicos", "Gestão de Estoques", 0)
)

$startRow = 204
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $r = $data[$i]
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    # keep the default "Normal" style on the new cells (column C otherwise
    # inherits the sheet's wide column style)
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 6)).Style = "Normal"
}
